$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repos column (B) for program families that were still missing one.
# Order matters for the shared-strings table append order: Schaap's repo
# first, then Shusheng's, then Keshav's (matches the authored commit).
$ws.Range("B10").Value = "https://github.com/aschaap/cas741.git"
$ws.Range("B6").Value = "https://github.com/sccdsyad8663/Shusheng-CAS741.git"
$ws.Range("B4").Value = "https://github.com/keshavd/cas741.git"

# Problem Statement Approved column (J) marked "Yes" for several rows.
$ws.Range("J6").Value = "Yes"
$ws.Range("J9").Value = "Yes"
$ws.Range("J11").Value = "Yes"
$ws.Range("J12").Value = "Yes"
$ws.Range("J13").Value = "Yes"

# Row-height tweaks around the newly edited rows.
$ws.Rows.Item(4).RowHeight = 13
$ws.Rows.Item(5).RowHeight = 17

# Move the active selection like the author's last cursor position.
$null = $ws.Range("J14").Select()
